$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(2, 8, 14, 20, 26, 32, 38, 44, 50, 56, 62, 68, 74)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value() -eq "2026/01/17") {
        $cell.NumberFormat = "@"
        $cell.Value = "2026/01/18"
        $cell.ClearFormats()
    }
}
